# Applies the cryptos.xlsx price-refresh edit described by the commit diff.
# Updates Price (D) / Volume(1h) (E) figures for most rows, and for a handful
# of rows (rank ties re-sorted upstream) also swaps the Coin (B) / Link (C)
# values between adjacent rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "80.650.45"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "3.119.07"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'204.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").Value = "'618.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").Value = "'0.277"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +22.25%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("D9").Value = "'0.573"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").Value = "3.128.65"
$ws.Range("E10").Value = "  -2.32%  "
$ws.Range("D11").Value = "'0.569"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").Value = "'0.0000247"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +9.32%  "
$ws.Range("D13").Value = "'0.164"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "'5.22"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.73%  "
$ws.Range("D15").Value = "3.708.15"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "'30.89"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "80.896.53"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").Value = "3.124.60"
$ws.Range("E18").Value = "  -2.59%  "
$ws.Range("D19").Value = "'3.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +9.87%  "
$ws.Range("D20").Value = "'13.77"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.72%  "
$ws.Range("D21").Value = "'426.28"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").Value = "'8.83"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -6.61%  "
$ws.Range("D23").Value = "'5.01"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +4.22%  "
$ws.Range("D25").Value = "'5.08"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.85%  "
$ws.Range("D26").Value = "3.305.36"
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("D27").Value = "'75.22"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("D28").Value = "'10.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.23%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "'0.0000118"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.27%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").Value = "'8.85"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").Value = "'547.31"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.91%  "
$ws.Range("D34").Value = "'1.45"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").Value = "'0.146"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +12.71%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.148"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.73%  "
$ws.Range("B37").Value = "PancakeSwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D37").Value = "'1.97"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").Value = "'22.36"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "'0.400"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").Value = "'20.69"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.42%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.81"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.87%  "
$ws.Range("D43").Value = "'2.96"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +16.91%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'160.26"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.45%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.96"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.85%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'184.35"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.89%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "'43.70"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").Value = "'1.30"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'0.761"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.65%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D51").Value = "'4.15"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.72%  "
